$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.949.97"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.211.14"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'289.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "'87.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'30.30"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'0.0774"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "'6.44"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "2.554.68"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "'13.91"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "2.207.03"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "39.901.95"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'11.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.28%  "
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'65.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'235.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").Value = "'22.43"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'9.17"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'155.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "'31.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'4.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "'0.0714"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  +6.34%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "'15.77"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").Value = "2.101.08"
$ws.Range("E41").Value = "  +7.69%  "
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").Value = "'9.96"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.79%  "
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "'17.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.81%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "2.428.42"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'68.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").Value = "'1.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.41%  "
